$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1242.3077
$ws.Range("I11").Value = 1242.3077
$ws.Range("K11").Value = 1242.3077
$ws.Range("M11").Value = -1102.3077
$ws.Range("H17").Value = 2812.5
$ws.Range("J17").Value = 2812.5
$ws.Range("L17").Value = 8437.5
$ws.Range("N17").Value = -8773.5
$ws.Range("H19").Value = 786.7917
$ws.Range("I19").Value = 346.6875
$ws.Range("J19").Value = 1667
$ws.Range("K19").Value = 346.6875
$ws.Range("L19").Value = 1667
$ws.Range("M19").Value = -171.6875
$ws.Range("N19").Value = -2017
$ws.Range("H28").Value = 2147.45
$ws.Range("I28").Value = 550.38464
$ws.Range("K28").Value = 550.38464
$ws.Range("M28").Value = -65.38463999999999
$ws.Range("H32").Value = 3040.1177
$ws.Range("I32").Value = 1422.7142
$ws.Range("J32").Value = 3459.4443
$ws.Range("K32").Value = 1422.7142
$ws.Range("L32").Value = 3459.4443
$ws.Range("M32").Value = -1096.7142
$ws.Range("N32").Value = -4111.4443
$ws.Range("H33").Value = 859.4545000000001
$ws.Range("J33").Value = 1324.5
$ws.Range("L33").Value = 1324.5
$ws.Range("N33").Value = -1782.5
$ws.Range("H40").Value = 4312
$ws.Range("I40").Value = 2000.3334
$ws.Range("J40").Value = 4658.75
$ws.Range("K40").Value = 2000.3334
$ws.Range("L40").Value = 4658.75
$ws.Range("M40").Value = -1825.3334
$ws.Range("N40").Value = -5008.75
$ws.Range("H43").Value = 1699.5834
$ws.Range("I43").Value = 1299.5714
$ws.Range("J43").Value = 2259.6
$ws.Range("K43").Value = 1299.5714
$ws.Range("L43").Value = 2259.6
$ws.Range("M43").Value = -1230.5714
$ws.Range("N43").Value = -2397.6
$ws.Range("H51").Value = 4817.173
$ws.Range("I51").Value = 2965.077
$ws.Range("J51").Value = 6669.269
$ws.Range("K51").Value = 2965.077
$ws.Range("L51").Value = 6669.269
$ws.Range("M51").Value = -2481.077
$ws.Range("N51").Value = -7637.269
$ws.Range("H55").Value = 204.52174
$ws.Range("I55").Value = 111.545456
$ws.Range("J55").Value = 289.75
$ws.Range("K55").Value = 111.545456
$ws.Range("L55").Value = 289.75
$ws.Range("M55").Value = 102.454544
$ws.Range("N55").Value = -717.75
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H62").Value = 8025.25
$ws.Range("I62").Value = 7929.923
$ws.Range("K62").Value = 7929.923
$ws.Range("M62").Value = -7305.923
$ws.Range("H64").Value = 8154.758
$ws.Range("J64").Value = 8309.936
$ws.Range("L64").Value = 8309.936
$ws.Range("N64").Value = -8805.936
$ws.Range("H65").Value = 8025.25
$ws.Range("I65").Value = 7929.923
$ws.Range("K65").Value = 39649.615
$ws.Range("M65").Value = -36529.615
$ws.Range("H67").Value = 8154.758
$ws.Range("J67").Value = 8309.936
$ws.Range("L67").Value = 8309.936
$ws.Range("N67").Value = -10025.936
$ws.Range("H92").Value = 975.2778
$ws.Range("I92").Value = 763.125
$ws.Range("J92").Value = 2672.5
$ws.Range("K92").Value = 763.125
$ws.Range("L92").Value = 2672.5
$ws.Range("M92").Value = 484.875
$ws.Range("N92").Value = -5168.5
$ws.Range("H97").Value = 2776.0588
$ws.Range("J97").Value = 2887.0625
$ws.Range("L97").Value = 8661.1875
$ws.Range("N97").Value = -9653.1875
$ws.Range("H103").Value = 737.0833
$ws.Range("I103").Value = 866.6667
$ws.Range("K103").Value = 2600.0001
$ws.Range("M103").Value = -2014.0001
$ws.Range("H116").Value = 4493.645
$ws.Range("I116").Value = 2934.3333
$ws.Range("J116").Value = 6652.6924
$ws.Range("K116").Value = 2934.3333
$ws.Range("L116").Value = 6652.6924
$ws.Range("M116").Value = 507.6667000000002
$ws.Range("N116").Value = -13536.6924
$ws.Range("H135").Value = 3225.6667
$ws.Range("I135").Value = 1799.3334
$ws.Range("J135").Value = 4176.5557
$ws.Range("K135").Value = 16194.0006
$ws.Range("L135").Value = 37589.0013
$ws.Range("M135").Value = -13659.0006
$ws.Range("N135").Value = -42659.0013
$ws.Range("H137").Value = 49244.668
$ws.Range("I137").Value = 68946.14999999999
$ws.Range("J137").Value = 4916.3335
$ws.Range("K137").Value = 206838.45
$ws.Range("L137").Value = 14749.0005
$ws.Range("M137").Value = -204288.45
$ws.Range("N137").Value = -19849.0005
$ws.Range("H138").Value = 4711.485
$ws.Range("I138").Value = 7600
$ws.Range("J138").Value = 4525.129
$ws.Range("K138").Value = 22800
$ws.Range("L138").Value = 13575.387
$ws.Range("M138").Value = -17660
$ws.Range("N138").Value = -23855.387
$ws.Range("H141").Value = 21593.8
$ws.Range("I141").Value = 25995.25
$ws.Range("J141").Value = 3988
$ws.Range("K141").Value = 77985.75
$ws.Range("L141").Value = 11964
$ws.Range("M141").Value = -72805.75
$ws.Range("N141").Value = -22324

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11020.087
$ws.Range("I32").Value = 8394.273999999999
$ws.Range("K32").Value = 8394.273999999999
$ws.Range("M32").Value = -8107.273999999999
$ws.Range("H61").Value = 4045.5588
$ws.Range("I61").Value = 3851.6667
$ws.Range("J61").Value = 5499.75
$ws.Range("K61").Value = 3851.6667
$ws.Range("L61").Value = 5499.75
$ws.Range("M61").Value = -3639.6667
$ws.Range("N61").Value = -5923.75
$ws.Range("H63").Value = 8382.916999999999
$ws.Range("I63").Value = 7339.2
$ws.Range("J63").Value = 9128.429
$ws.Range("K63").Value = 7339.2
$ws.Range("L63").Value = 9128.429
$ws.Range("M63").Value = -6653.2
$ws.Range("N63").Value = -10500.429
$ws.Range("H66").Value = 8382.916999999999
$ws.Range("I66").Value = 7339.2
$ws.Range("J66").Value = 9128.429
$ws.Range("K66").Value = 36696
$ws.Range("L66").Value = 45642.145
$ws.Range("M66").Value = -33264
$ws.Range("N66").Value = -52506.145
$ws.Range("H74").Value = 29451.705
$ws.Range("I74").Value = 3049.3333
$ws.Range("J74").Value = 59154.375
$ws.Range("K74").Value = 3049.3333
$ws.Range("L74").Value = 59154.375
$ws.Range("M74").Value = -2175.3333
$ws.Range("N74").Value = -60902.375
$ws.Range("H77").Value = 29451.705
$ws.Range("I77").Value = 3049.3333
$ws.Range("J77").Value = 59154.375
$ws.Range("K77").Value = 15246.6665
$ws.Range("L77").Value = 295771.875
$ws.Range("M77").Value = -10878.6665
$ws.Range("N77").Value = -304507.875
$ws.Range("H94").Value = 19999
$ws.Range("J94").Value = 19999
$ws.Range("L94").Value = 19999
$ws.Range("N94").Value = -21801
$ws.Range("H96").Value = 50000
$ws.Range("J96").Value = 50000
$ws.Range("L96").Value = 50000
$ws.Range("N96").Value = -55492
$ws.Range("H136").Value = 4045.5588
$ws.Range("I136").Value = 3851.6667
$ws.Range("J136").Value = 5499.75
$ws.Range("K136").Value = 11555.0001
$ws.Range("L136").Value = 16499.25
$ws.Range("M136").Value = -9005.000100000001
$ws.Range("N136").Value = -21599.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8405176
$ws.Range("I99").Value = 17858246
$ws.Range("K99").Value = 17858246
$ws.Range("M99").Value = -17856748
$ws.Range("H105").Value = 10419238
$ws.Range("I105").Value = 15627980
$ws.Range("K105").Value = 15627980
$ws.Range("M105").Value = -15626233
$ws.Range("H134").Value = 7344.8335
$ws.Range("I134").Value = 2428.05
$ws.Range("J134").Value = 17178.4
$ws.Range("K134").Value = 7284.150000000001
$ws.Range("L134").Value = 51535.2
$ws.Range("M134").Value = -4749.150000000001
$ws.Range("N134").Value = -56605.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1204.0667
$ws.Range("I16").Value = 871.375
$ws.Range("J16").Value = 1584.2858
$ws.Range("K16").Value = 871.375
$ws.Range("L16").Value = 1584.2858
$ws.Range("M16").Value = -584.375
$ws.Range("N16").Value = -2158.2858
$ws.Range("H58").Value = 6071.7144
$ws.Range("I58").Value = 6957.737
$ws.Range("K58").Value = 6957.737
$ws.Range("M58").Value = -6754.737
$ws.Range("H62").Value = 7818.5
$ws.Range("J62").Value = 7774.5
$ws.Range("L62").Value = 7774.5
$ws.Range("N62").Value = -9022.5
$ws.Range("H65").Value = 7818.5
$ws.Range("J65").Value = 7774.5
$ws.Range("L65").Value = 38872.5
$ws.Range("N65").Value = -45112.5
$ws.Range("H105").Value = 738.125
$ws.Range("I105").Value = 700.6667
$ws.Range("K105").Value = 700.6667
$ws.Range("M105").Value = 1046.3333
$ws.Range("H113").Value = 1204.0667
$ws.Range("I113").Value = 871.375
$ws.Range("J113").Value = 1584.2858
$ws.Range("K113").Value = 871.375
$ws.Range("L113").Value = 1584.2858
$ws.Range("M113").Value = 1298.625
$ws.Range("N113").Value = -5924.2858
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 103817.945
$ws.Range("I132").Value = 81205.30499999999
$ws.Range("K132").Value = 243615.915
$ws.Range("M132").Value = -241085.915
$ws.Range("H134").Value = 2834.3845
$ws.Range("I134").Value = 1777.2354
$ws.Range("J134").Value = 4831.222
$ws.Range("K134").Value = 5331.706200000001
$ws.Range("L134").Value = 14493.666
$ws.Range("M134").Value = -2796.706200000001
$ws.Range("N134").Value = -19563.666
$ws.Range("H135").Value = 70000
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H136").Value = 6071.7144
$ws.Range("I136").Value = 6957.737
$ws.Range("K136").Value = 20873.211
$ws.Range("M136").Value = -18323.211
$ws.Range("H138").Value = 27960.125
$ws.Range("J138").Value = 27960.125
$ws.Range("L138").Value = 27960.125
$ws.Range("N138").Value = -38240.125
$ws.Range("H141").Value = 390000.84
$ws.Range("J141").Value = 390000.84
$ws.Range("L141").Value = 390000.84
$ws.Range("N141").Value = -400360.84

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 54452.58
$ws.Range("J5").Value = 128125.25
$ws.Range("L5").Value = 384375.75
$ws.Range("N5").Value = -384599.75
$ws.Range("H12").Value = 35011.19
$ws.Range("J12").Value = 1878.3636
$ws.Range("L12").Value = 5635.0908
$ws.Range("N12").Value = -5981.0908
$ws.Range("H38").Value = 193.5
$ws.Range("I38").Value = 27
$ws.Range("J38").Value = 360
$ws.Range("K38").Value = 81
$ws.Range("L38").Value = 1080
$ws.Range("M38").Value = 266
$ws.Range("N38").Value = -1774
$ws.Range("H98").Value = 694.913
$ws.Range("J98").Value = 728.4666999999999
$ws.Range("L98").Value = 2185.4001
$ws.Range("N98").Value = -5181.4001
$ws.Range("H131").Value = 10423410
$ws.Range("I131").Value = 7578120.5
$ws.Range("J131").Value = 11502658
$ws.Range("K131").Value = 22734361.5
$ws.Range("L131").Value = 34507974
$ws.Range("M131").Value = -22729321.5
$ws.Range("N131").Value = -34518054
$ws.Range("H135").Value = 54452.58
$ws.Range("J135").Value = 128125.25
$ws.Range("L135").Value = 1153127.25
$ws.Range("N135").Value = -1158197.25
$ws.Range("H139").Value = 33334548
$ws.Range("I139").Value = 35715584
$ws.Range("J139").Value = 33
$ws.Range("K139").Value = 107146752
$ws.Range("L139").Value = 99
$ws.Range("M139").Value = -107141612
$ws.Range("N139").Value = -10379
$ws.Range("H141").Value = 4998.6665
$ws.Range("I141").Value = 4998.6665
$ws.Range("K141").Value = 14995.9995
$ws.Range("M141").Value = -9815.999500000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3000
$ws.Range("J23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3446
$ws.Range("H70").Value = 18195766
$ws.Range("I70").Value = 25004088
$ws.Range("J70").Value = 40236
$ws.Range("K70").Value = 25004088
$ws.Range("L70").Value = 40236
$ws.Range("M70").Value = -25003818
$ws.Range("N70").Value = -40776
$ws.Range("H73").Value = 18195766
$ws.Range("I73").Value = 25004088
$ws.Range("J73").Value = 40236
$ws.Range("K73").Value = 25004088
$ws.Range("L73").Value = 40236
$ws.Range("M73").Value = -25003152
$ws.Range("N73").Value = -42108
$ws.Range("H80").Value = 1529157.4
$ws.Range("I80").Value = 2136855
$ws.Range("K80").Value = 2136855
$ws.Range("M80").Value = -2135857
$ws.Range("H83").Value = 1529157.4
$ws.Range("I83").Value = 2136855
$ws.Range("K83").Value = 10684275
$ws.Range("M83").Value = -10679283
$ws.Range("H113").Value = 7986433.5
$ws.Range("J113").Value = 3277.5
$ws.Range("L113").Value = 3277.5
$ws.Range("N113").Value = -7617.5
$ws.Range("H122").Value = 358948.3
$ws.Range("I122").Value = 447440
$ws.Range("J122").Value = 4981.6
$ws.Range("K122").Value = 1342320
$ws.Range("L122").Value = 14944.8
$ws.Range("M122").Value = -1339870
$ws.Range("N122").Value = -19844.8
$ws.Range("H126").Value = 4157832.5
$ws.Range("I126").Value = 2676246
$ws.Range("J126").Value = 5956901.5
$ws.Range("K126").Value = 8028738
$ws.Range("L126").Value = 17870704.5
$ws.Range("M126").Value = -8026268
$ws.Range("N126").Value = -17875644.5
$ws.Range("H132").Value = 3945.182
$ws.Range("I132").Value = 3269.111
$ws.Range("J132").Value = 6987.5
$ws.Range("K132").Value = 9807.332999999999
$ws.Range("L132").Value = 20962.5
$ws.Range("M132").Value = -7277.332999999999
$ws.Range("N132").Value = -26022.5
$ws.Range("H134").Value = 52908
$ws.Range("J134").Value = 52908
$ws.Range("L134").Value = 158724
$ws.Range("N134").Value = -163794

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 70161.53999999999
$ws.Range("I22").Value = 112280.625
$ws.Range("K22").Value = 112280.625
$ws.Range("M22").Value = -111985.625
$ws.Range("H27").Value = 70161.53999999999
$ws.Range("I27").Value = 112280.625
$ws.Range("K27").Value = 112280.625
$ws.Range("M27").Value = -112173.625
$ws.Range("H46").Value = 6536.6
$ws.Range("I46").Value = 1350
$ws.Range("K46").Value = 1350
$ws.Range("M46").Value = -1162
$ws.Range("H55").Value = 2126.8235
$ws.Range("I55").Value = 2103.3333
$ws.Range("J55").Value = 2139.6365
$ws.Range("K55").Value = 2103.3333
$ws.Range("L55").Value = 2139.6365
$ws.Range("M55").Value = -1930.3333
$ws.Range("N55").Value = -2485.6365
$ws.Range("H61").Value = 7408155.5
$ws.Range("I61").Value = 11111877
$ws.Range("K61").Value = 11111877
$ws.Range("M61").Value = -11111675
$ws.Range("H93").Value = 41683170
$ws.Range("I93").Value = 55559390
$ws.Range("K93").Value = 55559390
$ws.Range("M93").Value = -55558142
$ws.Range("H100").Value = 4742.3335
$ws.Range("I100").Value = 4900
$ws.Range("K100").Value = 4900
$ws.Range("M100").Value = -4359
$ws.Range("H104").Value = 39037
$ws.Range("J104").Value = 39037
$ws.Range("L104").Value = 39037
$ws.Range("N104").Value = -46025
$ws.Range("H113").Value = 7408155.5
$ws.Range("I113").Value = 11111877
$ws.Range("K113").Value = 11111877
$ws.Range("M113").Value = -11109707
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 6811
$ws.Range("I132").Value = 7183.3193
$ws.Range("K132").Value = 21549.9579
$ws.Range("M132").Value = -19019.9579
$ws.Range("H136").Value = 65859.03
$ws.Range("I136").Value = 74447.82000000001
$ws.Range("K136").Value = 223343.46
$ws.Range("M136").Value = -220793.46

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 50000000
$ws.Range("I5").Value = 50000000
$ws.Range("K5").Value = 50000000
$ws.Range("M5").Value = -49999888
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H46").Value = 74344.5
$ws.Range("I46").Value = 59999
$ws.Range("J46").Value = 79126.336
$ws.Range("K46").Value = 59999
$ws.Range("L46").Value = 79126.336
$ws.Range("M46").Value = -59768
$ws.Range("N46").Value = -79588.336
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H100").Value = 6678.143
$ws.Range("I100").Value = 7641.3335
$ws.Range("K100").Value = 15282.667
$ws.Range("M100").Value = -14741.667
$ws.Range("H113").Value = 1316.8334
$ws.Range("I113").Value = 631.5
$ws.Range("J113").Value = 2173.5
$ws.Range("K113").Value = 1894.5
$ws.Range("L113").Value = 6520.5
$ws.Range("M113").Value = 275.5
$ws.Range("N113").Value = -10860.5
$ws.Range("H132").Value = 19825598
$ws.Range("I132").Value = 25004456
$ws.Range("J132").Value = 993389.0600000001
$ws.Range("K132").Value = 75013368
$ws.Range("L132").Value = 2980167.18
$ws.Range("M132").Value = -75010838
$ws.Range("N132").Value = -2985227.18
$ws.Range("H134").Value = 74344.5
$ws.Range("I134").Value = 59999
$ws.Range("J134").Value = 79126.336
$ws.Range("K134").Value = 179997
$ws.Range("L134").Value = 237379.008
$ws.Range("M134").Value = -177462
$ws.Range("N134").Value = -242449.008
$ws.Range("H136").Value = 2751.875
$ws.Range("I136").Value = 2453.4443
$ws.Range("J136").Value = 3972.7273
$ws.Range("K136").Value = 7360.3329
$ws.Range("L136").Value = 11918.1819
$ws.Range("M136").Value = -4810.3329
$ws.Range("N136").Value = -17018.1819
